$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first worksheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 170
$wsExhibit.Range("F4").Value = 754

# Sheet "全部类型" (All types) - fourth worksheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 170
$wsAll.Range("F5").Value = 754
